# Auto-generated: append SQL schema section after the trailing empty paragraph
$d = $word.ActiveDocument
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>SQL:</w:t>
      </w:r>
      <w:r>
        <w:br/>
      </w:r>
      <w:r>
        <w:t>CREATE TABLE airline (</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>airline_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> SERIAL PRIMARY KEY,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>carrier_name</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>VARCHAR(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>50)</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t>);</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t>CREATE TABLE airport (</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>airport_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> SERIAL PRIMARY KEY,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    airport </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>VARCHAR(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>5),</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    city </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>VARCHAR(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>50)</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t>);</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:lastRenderedPageBreak/>
        <w:t>CREATE TABLE flight (</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>flight_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> SERIAL PRIMARY KEY,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    month INT,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>total_arrivals</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> FLOAT,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>total_delays_ct</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> FLOAT,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>flight_cancelled</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> FLOAT,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>flight_diverted</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> FLOAT,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>total_delays_min</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> FLOAT,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>airline_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> INT REFERENCES airline(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>airline_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>),</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>airport_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> INT REFERENCES airport(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>airport_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>)</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t>);</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t>CREATE TABLE delay (</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>delay_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> SERIAL PRIMARY KEY,</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>flight_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> INT REFERENCES flight(</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>flight_id</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>),</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>delay_type</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>VARCHAR(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>50),</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t xml:space="preserve">    </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>delay_duration</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> FLOAT</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

$r.InsertParagraphAfter()
$r = $d.Paragraphs.Last.Range
$r.Collapse(0)
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r>
        <w:t>);</w:t>
      </w:r></w:p>
'@
$r.InsertXML($xml)
$r = $d.Paragraphs.Last.Range

Write-Output ("paragraphs=" + $d.Paragraphs.Count)
